# Applies the changes described by the diff:
#  - C11: 5 -> 9
#  - C18: (empty) -> 5
#  - C19: (empty) -> 10
#  - C51 (=SUM(C6:C50)) recalculates to 64
#  - sheet view: scrolled/selected cell moved (topLeftCell A4->A7, selection C22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = 9
$ws.Range("C18").Value = 5
$ws.Range("C19").Value = 10

# Update the visible scroll position / active selection to match the
# author's saved view state.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C22").Select()
